$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.799.50'
$ws.Range("E2").Value = '  -0.55%  '
$ws.Range("D3").Value = '2.947.41'
$ws.Range("E3").Value = '  +2.64%  '
$ws.Range("E4").Value = '  +0.15%  '
$ws.Range("D5").Value = '''352.21'
$ws.Range("E5").Value = '  +0.58%  '
$ws.Range("D6").Value = '''111.57'
$ws.Range("E6").Value = '  -0.74%  '
$ws.Range("D7").Value = '''0.563'
$ws.Range("E7").Value = '  +0.80%  '
$ws.Range("E8").Value = '  +0.08%  '
$ws.Range("E9").Value = '  +2.36%  '
$ws.Range("D10").Value = '''39.48'
$ws.Range("E10").Value = '  -1.78%  '
$ws.Range("D11").Value = '''0.0897'
$ws.Range("E11").Value = '  +5.50%  '
$ws.Range("E12").Value = '  +0.89%  '
$ws.Range("D13").Value = '''19.84'
$ws.Range("E13").Value = '  -1.03%  '
$ws.Range("D14").Value = '''8.04'
$ws.Range("E14").Value = '  +2.49%  '
$ws.Range("D15").Value = '3.416.25'
$ws.Range("E15").Value = '  +2.74%  '
$ws.Range("D16").Value = '2.956.82'
$ws.Range("E16").Value = '  +2.67%  '
$ws.Range("D17").Value = '''0.996'
$ws.Range("E17").Value = '  +0.21%  '
$ws.Range("D18").Value = '51.933.36'
$ws.Range("E18").Value = '  -0.26%  '
$ws.Range("D19").Value = '''7.70'
$ws.Range("E19").Value = '  +0.84%  '
$ws.Range("D20").Value = '''14.54'
$ws.Range("E20").Value = '  +7.01%  '
$ws.Range("E21").Value = '  -2.45%  '
$ws.Range("D22").Value = '0.0₃0988'
$ws.Range("E22").Value = '  +1.55%  '
$ws.Range("D23").Value = '''71.42'
$ws.Range("E23").Value = '  +0.76%  '
$ws.Range("D24").Value = '''273.07'
$ws.Range("E24").Value = '  +1.06%  '
$ws.Range("E25").Value = '  +0.63%  '
$ws.Range("E26").Value = '  +11.92%  '
$ws.Range("D27").Value = '''27.37'
$ws.Range("E27").Value = '  +3.23%  '
$ws.Range("E28").Value = '  +0.05%  '
$ws.Range("D29").Value = '''7.41'
$ws.Range("E29").Value = '  +18.46%  '
$ws.Range("E30").Value = '  +23.46%  '
$ws.Range("D31").Value = '''10.76'
$ws.Range("E31").Value = '  +1.83%  '
$ws.Range("D32").Value = '''6.35'
$ws.Range("E32").Value = '  +8.78%  '
$ws.Range("D33").Value = '''37.56'
$ws.Range("E33").Value = '  -3.22%  '
$ws.Range("D34").Value = '''53.09'
$ws.Range("E34").Value = '  +1.42%  '
$ws.Range("E35").Value = '  -1.16%  '
$ws.Range("E36").Value = '  -0.07%  '
$ws.Range("B37").Value = 'Toncoin'
$ws.Range("C37").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D37").Value = '''1.87'
$ws.Range("E37").Value = '  -1.06%  '
$ws.Range("B38").Value = 'LidoDAOToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D38").Value = '''3.39'
$ws.Range("E38").Value = '  +2.23%  '
$ws.Range("D39").Value = '''18.78'
$ws.Range("E39").Value = '  +0.57%  '
$ws.Range("E40").Value = '  +1.11%  '
$ws.Range("E41").Value = '  +0.31%  '
$ws.Range("E42").Value = '  +2.43%  '
$ws.Range("D43").Value = '''23.54'
$ws.Range("E43").Value = '  +4.95%  '
$ws.Range("E44").Value = '  -1.56%  '
$ws.Range("D45").Value = '''3.54'
$ws.Range("E45").Value = '  +0.95%  '
$ws.Range("E46").Value = '  +2.07%  '
$ws.Range("D47").Value = '2.162.92'
$ws.Range("E47").Value = '  -0.37%  '
$ws.Range("D48").Value = '''113.41'
$ws.Range("E48").Value = '  -7.03%  '
$ws.Range("E49").Value = '  +2.89%  '
$ws.Range("E50").Value = '  +4.68%  '
$ws.Range("D51").Value = '''0.924'
$ws.Range("E51").Value = '  -3.84%  '
